# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Update DAMSLTag (column I) and DialogAct
# (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Range("I20").Value = "ba"
$ws.Range("J20").Value = "Appreciation"

# Row 53
$ws.Range("I53").Value = "ba"
$ws.Range("J53").Value = "Appreciation"

# Row 54
$ws.Range("I54").Value = "sd"
$ws.Range("J54").Value = "Statement-non-opinion"

# Row 61
$ws.Range("I61").Value = "%"
$ws.Range("J61").Value = "Uninterpretable"

# Row 69
$ws.Range("I69").Value = "%"
$ws.Range("J69").Value = "Uninterpretable"

# Row 82
$ws.Range("I82").Value = "%"
$ws.Range("J82").Value = "Uninterpretable"
